# Applies the "cryptos list" refresh captured in the commit:
#   "Updated cryptos list on Sun Jul 28 16:30:30 UTC 2024 with GitHub Actions"
#
# Each entry updates one cell (Row/Col, 1-based, matching Cells.Item(row, col)).
# AsText cells hold numeric-looking strings (e.g. "67.677.30", "6.58") that must
# stay plain text (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel, so they are written via Formula with a
# leading text-prefix apostrophe (') rather than via Value.
$updates = @(
    @{ Row = 2; Col = 4; Cell = "D2"; Value = "67.677.30"; AsText = $true }
    @{ Row = 2; Col = 5; Cell = "E2"; Value = "  -1.76%  "; AsText = $false }
    @{ Row = 3; Col = 4; Cell = "D3"; Value = "3.270.69"; AsText = $true }
    @{ Row = 3; Col = 5; Cell = "E3"; Value = "  -0.34%  "; AsText = $false }
    @{ Row = 4; Col = 5; Cell = "E4"; Value = "  +0.02%  "; AsText = $false }
    @{ Row = 5; Col = 4; Cell = "D5"; Value = "580.96"; AsText = $true }
    @{ Row = 5; Col = 5; Cell = "E5"; Value = "  -0.78%  "; AsText = $false }
    @{ Row = 6; Col = 4; Cell = "D6"; Value = "184.53"; AsText = $true }
    @{ Row = 6; Col = 5; Cell = "E6"; Value = "  +1.02%  "; AsText = $false }
    @{ Row = 7; Col = 5; Cell = "E7"; Value = "  +0.03%  "; AsText = $false }
    @{ Row = 8; Col = 4; Cell = "D8"; Value = "0.604"; AsText = $true }
    @{ Row = 8; Col = 5; Cell = "E8"; Value = "  +0.83%  "; AsText = $false }
    @{ Row = 9; Col = 5; Cell = "E9"; Value = "  -3.32%  "; AsText = $false }
    @{ Row = 10; Col = 4; Cell = "D10"; Value = "6.58"; AsText = $true }
    @{ Row = 10; Col = 5; Cell = "E10"; Value = "  -1.37%  "; AsText = $false }
    @{ Row = 11; Col = 5; Cell = "E11"; Value = "  -3.84%  "; AsText = $false }
    @{ Row = 12; Col = 4; Cell = "D12"; Value = "3.836.40"; AsText = $true }
    @{ Row = 12; Col = 5; Cell = "E12"; Value = "  -0.34%  "; AsText = $false }
    @{ Row = 13; Col = 4; Cell = "D13"; Value = "0.138"; AsText = $true }
    @{ Row = 13; Col = 5; Cell = "E13"; Value = "  +0.84%  "; AsText = $false }
    @{ Row = 14; Col = 4; Cell = "D14"; Value = "27.34"; AsText = $true }
    @{ Row = 14; Col = 5; Cell = "E14"; Value = "  -4.97%  "; AsText = $false }
    @{ Row = 15; Col = 4; Cell = "D15"; Value = "67.740.08"; AsText = $true }
    @{ Row = 15; Col = 5; Cell = "E15"; Value = "  -1.59%  "; AsText = $false }
    @{ Row = 16; Col = 5; Cell = "E16"; Value = "  -2.56%  "; AsText = $false }
    @{ Row = 17; Col = 4; Cell = "D17"; Value = "3.266.01"; AsText = $true }
    @{ Row = 17; Col = 5; Cell = "E17"; Value = "  -0.23%  "; AsText = $false }
    @{ Row = 18; Col = 5; Cell = "E18"; Value = "  -2.27%  "; AsText = $false }
    @{ Row = 19; Col = 5; Cell = "E19"; Value = "  -1.57%  "; AsText = $false }
    @{ Row = 20; Col = 4; Cell = "D20"; Value = "403.05"; AsText = $true }
    @{ Row = 20; Col = 5; Cell = "E20"; Value = "  +1.50%  "; AsText = $false }
    @{ Row = 21; Col = 4; Cell = "D21"; Value = "7.55"; AsText = $true }
    @{ Row = 21; Col = 5; Cell = "E21"; Value = "  -2.38%  "; AsText = $false }
    @{ Row = 22; Col = 5; Cell = "E22"; Value = "  +0.12%  "; AsText = $false }
    @{ Row = 23; Col = 4; Cell = "D23"; Value = "70.89"; AsText = $true }
    @{ Row = 23; Col = 5; Cell = "E23"; Value = "  -1.29%  "; AsText = $false }
    @{ Row = 24; Col = 5; Cell = "E24"; Value = "  -1.69%  "; AsText = $false }
    @{ Row = 25; Col = 5; Cell = "E25"; Value = "  -2.54%  "; AsText = $false }
    @{ Row = 26; Col = 4; Cell = "D26"; Value = "0.186"; AsText = $true }
    @{ Row = 26; Col = 5; Cell = "E26"; Value = "  -1.66%  "; AsText = $false }
    @{ Row = 27; Col = 4; Cell = "D27"; Value = "9.51"; AsText = $true }
    @{ Row = 27; Col = 5; Cell = "E27"; Value = "  -2.24%  "; AsText = $false }
    @{ Row = 28; Col = 5; Cell = "E28"; Value = "  +0.54%  "; AsText = $false }
    @{ Row = 29; Col = 4; Cell = "D29"; Value = "1.94"; AsText = $true }
    @{ Row = 29; Col = 5; Cell = "E29"; Value = "  -2.10%  "; AsText = $false }
    @{ Row = 30; Col = 4; Cell = "D30"; Value = "22.61"; AsText = $true }
    @{ Row = 30; Col = 5; Cell = "E30"; Value = "  -1.79%  "; AsText = $false }
    @{ Row = 31; Col = 5; Cell = "E31"; Value = "  -4.67%  "; AsText = $false }
    @{ Row = 32; Col = 4; Cell = "D32"; Value = "6.91"; AsText = $true }
    @{ Row = 32; Col = 5; Cell = "E32"; Value = "  -3.44%  "; AsText = $false }
    @{ Row = 33; Col = 5; Cell = "E33"; Value = "  +0.06%  "; AsText = $false }
    @{ Row = 34; Col = 5; Cell = "E34"; Value = "  -4.23%  "; AsText = $false }
    @{ Row = 35; Col = 4; Cell = "D35"; Value = "163.66"; AsText = $true }
    @{ Row = 35; Col = 5; Cell = "E35"; Value = "  -0.67%  "; AsText = $false }
    @{ Row = 36; Col = 4; Cell = "D36"; Value = "1.46"; AsText = $true }
    @{ Row = 36; Col = 5; Cell = "E36"; Value = "  -3.65%  "; AsText = $false }
    @{ Row = 37; Col = 5; Cell = "E37"; Value = "  -1.24%  "; AsText = $false }
    @{ Row = 38; Col = 4; Cell = "D38"; Value = "26.91"; AsText = $true }
    @{ Row = 38; Col = 5; Cell = "E38"; Value = "  +2.02%  "; AsText = $false }
    @{ Row = 39; Col = 5; Cell = "E39"; Value = "  -3.08%  "; AsText = $false }
    @{ Row = 40; Col = 5; Cell = "E40"; Value = "  -1.96%  "; AsText = $false }
    @{ Row = 41; Col = 5; Cell = "E41"; Value = "  -3.75%  "; AsText = $false }
    @{ Row = 42; Col = 4; Cell = "D42"; Value = "2.671.69"; AsText = $true }
    @{ Row = 43; Col = 4; Cell = "D43"; Value = "40.71"; AsText = $true }
    @{ Row = 43; Col = 5; Cell = "E43"; Value = "  -1.83%  "; AsText = $false }
    @{ Row = 44; Col = 2; Cell = "B44"; Value = "dogwifhat"; AsText = $false }
    @{ Row = 44; Col = 3; Cell = "C44"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; AsText = $false }
    @{ Row = 44; Col = 4; Cell = "D44"; Value = "2.44"; AsText = $true }
    @{ Row = 44; Col = 5; Cell = "E44"; Value = "  -4.75%  "; AsText = $false }
    @{ Row = 45; Col = 2; Cell = "B45"; Value = "Hedera"; AsText = $false }
    @{ Row = 45; Col = 3; Cell = "C45"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; AsText = $false }
    @{ Row = 45; Col = 4; Cell = "D45"; Value = "0.0677"; AsText = $true }
    @{ Row = 45; Col = 5; Cell = "E45"; Value = "  -1.86%  "; AsText = $false }
    @{ Row = 46; Col = 2; Cell = "B46"; Value = "InjectiveProtocol"; AsText = $false }
    @{ Row = 46; Col = 3; Cell = "C46"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; AsText = $false }
    @{ Row = 46; Col = 4; Cell = "D46"; Value = "24.71"; AsText = $true }
    @{ Row = 46; Col = 5; Cell = "E46"; Value = "  -0.08%  "; AsText = $false }
    @{ Row = 47; Col = 2; Cell = "B47"; Value = "Bittensor"; AsText = $false }
    @{ Row = 47; Col = 3; Cell = "C47"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; AsText = $false }
    @{ Row = 47; Col = 4; Cell = "D47"; Value = "334.70"; AsText = $true }
    @{ Row = 47; Col = 5; Cell = "E47"; Value = "  -3.41%  "; AsText = $false }
    @{ Row = 48; Col = 4; Cell = "D48"; Value = "0.0274"; AsText = $true }
    @{ Row = 48; Col = 5; Cell = "E48"; Value = "  -3.17%  "; AsText = $false }
    @{ Row = 49; Col = 4; Cell = "D49"; Value = "6.27"; AsText = $true }
    @{ Row = 49; Col = 5; Cell = "E49"; Value = "  -0.74%  "; AsText = $false }
    @{ Row = 50; Col = 4; Cell = "D50"; Value = "0.100"; AsText = $true }
    @{ Row = 50; Col = 5; Cell = "E50"; Value = "  -1.81%  "; AsText = $false }
    @{ Row = 51; Col = 5; Cell = "E51"; Value = "  -1.88%  "; AsText = $false }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.AsText) {
        # Leading apostrophe forces Excel to store/keep the value as text.
        $cell.Formula = "'" + $u.Value
    } else {
        $cell.Value = $u.Value
    }
}

Write-Host "Applied $($updates.Count) cell updates"
